$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A132").Value = '```sql
SELECT "order_status" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" >= ''2025-07-01'' AND "order_created_at" < ''2025-08-01'' LIMIT 1;
```'
$ws.Rows(132).AutoFit()
$ws.Range("A133").Value = 'SELECT "order_status" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "final_status_at" >= ''2025-07-01 00:00:00'' AND "final_status_at" < ''2025-08-01 00:00:00'' ORDER BY "final_status_at" DESC LIMIT 1;'
$ws.Range("A134").Value = 'SELECT "id" AS order_number FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND EXTRACT(MONTH FROM "order_created_at") = 7 ORDER BY "order_created_at" DESC LIMIT 1;'
$ws.Range("A135").Value = 'SELECT "order_status" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "id" = 1823383 LIMIT 1;'
$ws.Range("A136").Value = 'SELECT "id", "order_status", "delivery_date", "order_created_at", "order_accepted_at", "start_ride_at", "reached_shop_at", "order_picked_at", "shipped_at", "reached_dest_at", "final_status_at", "customer_number", "customer_name", "client_id", "captain_name", "shop_to_delivery_km", "cancellation_reason" FROM "updated_table" WHERE "client_name" = ''MC DONALDS'' AND "order_created_at" >= ''2025-07-01'' AND "order_created_at" < ''2025-08-01'' LIMIT 1;'
